$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")
$ws.Activate()

# Insert a new column D ("ReadPopulationFromCSV") before the current ModelParameterSheets
# column, shifting ModelParameterSheets..OutputPathsIds one column to the right (D->E ... L->M).
$ws.Columns.Item(4).EntireColumn.Insert()

# Header for the new column (bold, like the rest of the header row).
$ws.Range("D1").Value = "ReadPopulationFromCSV"
$ws.Range("D1").Font.Bold = $true

# Row 2 (TestScenario): no ReadPopulationFromCSV value; the old SteadyState value
# (shifted from H2 to I2 by the column insert) is removed entirely for this scenario.
$ws.Range("I2").ClearContents()

# Row 4 (PopulationScenario): new ReadPopulationFromCSV value = FALSE (created from defined
# population demographics, the default behavior).
$ws.Range("D4").Value = $false

# Row 5 (new scenario): PopulationScenarioFromCSV - same as PopulationScenario but reads the
# population from a CSV file instead (ReadPopulationFromCSV = TRUE).
$ws.Range("A5").Value = "PopulationScenarioFromCSV"
$ws.Range("B5").Value = "Indiv"
$ws.Range("C5").Value = "TestPopulation"
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "Global"
$ws.Range("F5").Value = "Aciclovir_iv_250mg"
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = "h"
$ws.Range("I5").Value = $false
$ws.Range("L5").Value = "Aciclovir.pkml"

# Re-fit the column widths for the changed columns to reflect the new (longer) content.
$ws.Columns.Item(1).ColumnWidth = 17.6953125
$ws.Columns.Item(2).ColumnWidth = 10.6640625
$ws.Columns.Item(3).ColumnWidth = 13.59375
$ws.Columns.Item(4).ColumnWidth = 22.734375
$ws.Columns.Item(9).ColumnWidth = 10.6640625
$ws.Columns.Item(10).ColumnWidth = 15.3515625
$ws.Columns.Item(11).ColumnWidth = 19.3359375
$ws.Columns.Item(12).ColumnWidth = 12.890625
$ws.Columns.Item(13).ColumnWidth = 29.1796875

# Match the final selection left behind in the workbook.
[void]$ws.Range("G13").Select()
